$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets ---
$wsModel   = $wb.Worksheets.Item("ModelMetrics")
$wsFeature = $wb.Worksheets.Item("FeatureImportance")
$wsMeta    = $wb.Worksheets.Item("Metadata")

$wsModel.Name   = "Model Performance"
$wsFeature.Name = "Feature Importance"
$wsMeta.Name    = "Deployment Info"

# --- 2. Replace the contents of the (renamed) Deployment Info sheet ---
# Clear the old Field/Value metadata table entirely.
$wsMeta.Cells.Clear()

# New header row - deployment metadata for the Power BI export.
$wsMeta.Range("A1").Value = "deployment_date"
$wsMeta.Range("B1").Value = "model_version"
$wsMeta.Range("C1").Value = "data_source"
$wsMeta.Range("D1").Value = "refresh_frequency"

# New data row.
$wsMeta.Range("A2").Value = (Get-Date -Year 2025 -Month 10 -Day 14)
$wsMeta.Range("A2").NumberFormat = "mm-dd-yyyy"
$wsMeta.Range("B2").Value = "v1.0"
$wsMeta.Range("C2").Value = "HDPSA Clean Dataset"
$wsMeta.Range("D2").Value = "Weekly"
